# Auto-generated edit script applying numeric corrections to multiple sheets
# as described in the commit diff (scheduled runner value updates).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 19379.455  # H33: was 23645.889
$ws.Cells.Item(33, 9).Value = 30148.143  # I33: was 42135.2
$ws.Cells.Item(33, 11).Value = 30148.143  # K33: was 42135.2
$ws.Cells.Item(33, 13).Value = -29919.143  # M33: was -41906.2
$ws.Cells.Item(53, 8).Value = 1778.4546  # H53: was 1766.9
$ws.Cells.Item(53, 9).Value = 1462.3334  # I53: was 1479
$ws.Cells.Item(53, 10).Value = 2157.8  # J53: was 2198.75
$ws.Cells.Item(53, 11).Value = 1462.3334  # K53: was 1479
$ws.Cells.Item(53, 12).Value = 2157.8  # L53: was 2198.75
$ws.Cells.Item(53, 13).Value = -825.3334  # M53: was -842
$ws.Cells.Item(53, 14).Value = -3431.8  # N53: was -3472.75
$ws.Cells.Item(97, 8).Value = 1385.0714  # H97: was 1169.2354
$ws.Cells.Item(97, 10).Value = 1385.0714  # J97: was 1169.2354
$ws.Cells.Item(97, 12).Value = 4155.2142  # L97: was 3507.7062
$ws.Cells.Item(97, 14).Value = -5147.2142  # N97: was -4499.706200000001
$ws.Cells.Item(107, 8).Value = 431.51852  # H107: was 446.15384
$ws.Cells.Item(107, 9).Value = 545.55554  # I107: was 523.8421
$ws.Cells.Item(107, 10).Value = 203.44444  # J107: was 235.28572
$ws.Cells.Item(107, 11).Value = 545.55554  # K107: was 523.8421
$ws.Cells.Item(107, 12).Value = 203.44444  # L107: was 235.28572
$ws.Cells.Item(107, 13).Value = 1374.44446  # M107: was 1396.1579
$ws.Cells.Item(107, 14).Value = -4043.44444  # N107: was -4075.28572
$ws.Cells.Item(111, 8).Value = 2783.65  # H111: was 2878.0527
$ws.Cells.Item(111, 10).Value = 2329.6667  # J111: was 2999.5
$ws.Cells.Item(111, 12).Value = 6989.000100000001  # L111: was 8998.5
$ws.Cells.Item(111, 14).Value = -13123.0001  # N111: was -15132.5
$ws.Cells.Item(115, 8).Value = 748.1  # H115: was 808.1
$ws.Cells.Item(115, 9).Value = 720.1111  # I115: was 766.375
$ws.Cells.Item(115, 10).Value = 1000  # J115: was 975
$ws.Cells.Item(115, 11).Value = 2160.3333  # K115: was 2299.125
$ws.Cells.Item(115, 12).Value = 3000  # L115: was 2925
$ws.Cells.Item(115, 13).Value = -593.3332999999998  # M115: was -732.125
$ws.Cells.Item(115, 14).Value = -6134  # N115: was -6059
$ws.Cells.Item(132, 8).Value = 1915  # H132: was 1915.2094
$ws.Cells.Item(132, 9).Value = 1808.561  # I132: was 1808.7805
$ws.Cells.Item(132, 11).Value = 5425.683  # K132: was 5426.3415
$ws.Cells.Item(132, 13).Value = -2895.683  # M132: was -2896.3415
$ws.Cells.Item(135, 8).Value = 714.8333  # H135: was 727.25714
$ws.Cells.Item(135, 9).Value = 482.3846  # I135: was 490.48
$ws.Cells.Item(135, 11).Value = 4341.4614  # K135: was 4414.32
$ws.Cells.Item(135, 13).Value = -1806.4614  # M135: was -1879.32
$ws.Cells.Item(141, 8).Value = 2021.6  # H141: was 2023.1428
$ws.Cells.Item(141, 9).Value = 2165  # I141: was 2177.6924
$ws.Cells.Item(141, 11).Value = 6495  # K141: was 6533.0772
$ws.Cells.Item(141, 13).Value = -1315  # M141: was -1353.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5277.9604  # H32: was 5817.026
$ws.Cells.Item(32, 9).Value = 4384.2817  # I32: was 4601.2427
$ws.Cells.Item(32, 10).Value = 17968.2  # J32: was 16455.125
$ws.Cells.Item(32, 11).Value = 4384.2817  # K32: was 4601.2427
$ws.Cells.Item(32, 12).Value = 17968.2  # L32: was 16455.125
$ws.Cells.Item(32, 13).Value = -4097.2817  # M32: was -4314.2427
$ws.Cells.Item(32, 14).Value = -18542.2  # N32: was -17029.125
$ws.Cells.Item(45, 8).Value = 2079.4285  # H45: was 2641
$ws.Cells.Item(45, 9).Value = 1998.25  # I45: was 3500
$ws.Cells.Item(45, 10).Value = 2187.6667  # J45: was 2354.6667
$ws.Cells.Item(45, 11).Value = 1998.25  # K45: was 3500
$ws.Cells.Item(45, 12).Value = 2187.6667  # L45: was 2354.6667
$ws.Cells.Item(45, 13).Value = -1621.25  # M45: was -3123
$ws.Cells.Item(45, 14).Value = -2941.6667  # N45: was -3108.6667
$ws.Cells.Item(74, 8).Value = 34522964  # H74: was 33372258
$ws.Cells.Item(74, 9).Value = 50057284  # I74: was 45506704
$ws.Cells.Item(74, 10).Value = 2251.889  # J74: was 2538.375
$ws.Cells.Item(74, 11).Value = 50057284  # K74: was 45506704
$ws.Cells.Item(74, 12).Value = 2251.889  # L74: was 2538.375
$ws.Cells.Item(74, 13).Value = -50056410  # M74: was -45505830
$ws.Cells.Item(74, 14).Value = -3999.889  # N74: was -4286.375
$ws.Cells.Item(77, 8).Value = 34522964  # H77: was 33372258
$ws.Cells.Item(77, 9).Value = 50057284  # I77: was 45506704
$ws.Cells.Item(77, 10).Value = 2251.889  # J77: was 2538.375
$ws.Cells.Item(77, 11).Value = 250286420  # K77: was 227533520
$ws.Cells.Item(77, 12).Value = 11259.445  # L77: was 12691.875
$ws.Cells.Item(77, 13).Value = -250282052  # M77: was -227529152
$ws.Cells.Item(77, 14).Value = -19995.445  # N77: was -21427.875
$ws.Cells.Item(104, 8).Value = 81408  # H104: was 20000
$ws.Cells.Item(104, 10).Value = 81408  # J104: was 20000
$ws.Cells.Item(104, 12).Value = 81408  # L104: was 20000
$ws.Cells.Item(104, 14).Value = -88396  # N104: was -26988
$ws.Cells.Item(132, 8).Value = 58920000  # H132: was 66775684
$ws.Cells.Item(132, 9).Value = 16001.4  # I132: was 18094.924
$ws.Cells.Item(132, 11).Value = 48004.2  # K132: was 54284.772
$ws.Cells.Item(132, 13).Value = -45474.2  # M132: was -51754.772
$ws.Cells.Item(140, 8).Value = 150000  # H140: was 0
$ws.Cells.Item(140, 10).Value = 150000  # J140: was 0
$ws.Cells.Item(140, 12).Value = 150000  # L140: was 0
$ws.Cells.Item(140, 14).Value = -160360  # N140: new cell

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2474.625  # H20: was 2411.6667
$ws.Cells.Item(20, 9).Value = 2949.75  # I20: was 2741.4
$ws.Cells.Item(20, 11).Value = 2949.75  # K20: was 2741.4
$ws.Cells.Item(20, 13).Value = -2702.75  # M20: was -2494.4
$ws.Cells.Item(53, 8).Value = 0  # H53: was 55000
$ws.Cells.Item(53, 10).Value = 0  # J53: was 55000
$ws.Cells.Item(140, 8).Value = 159233.2  # H140: was 133701.2
$ws.Cells.Item(140, 10).Value = 159233.2  # J140: was 133701.2
$ws.Cells.Item(140, 12).Value = 159233.2  # L140: was 133701.2
$ws.Cells.Item(140, 14).Value = -169593.2  # N140: was -144061.2
$ws.Cells.Item(53, 12).Value = 0  # L53: was 55000
$ws.Cells.Item(53, 14).Value = $null  # N53: was -56148

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 217.52174  # H7: was 211.1923
$ws.Cells.Item(7, 9).Value = 68.083336  # I7: was 71.92308
$ws.Cells.Item(7, 10).Value = 380.54544  # J7: was 350.46155
$ws.Cells.Item(7, 11).Value = 68.083336  # K7: was 71.92308
$ws.Cells.Item(7, 12).Value = 380.54544  # L7: was 350.46155
$ws.Cells.Item(7, 13).Value = 44.916664  # M7: was 41.07692
$ws.Cells.Item(7, 14).Value = -606.54544  # N7: was -576.46155
$ws.Cells.Item(31, 8).Value = 3582.628  # H31: was 3632.2144
$ws.Cells.Item(31, 9).Value = 2690.5652  # I31: was 2744.682
$ws.Cells.Item(31, 11).Value = 2690.5652  # K31: was 2744.682
$ws.Cells.Item(31, 13).Value = -2395.5652  # M31: was -2449.682
$ws.Cells.Item(34, 8).Value = 3582.628  # H34: was 3632.2144
$ws.Cells.Item(34, 9).Value = 2690.5652  # I34: was 2744.682
$ws.Cells.Item(34, 11).Value = 2690.5652  # K34: was 2744.682
$ws.Cells.Item(34, 13).Value = -2488.5652  # M34: was -2542.682
$ws.Cells.Item(132, 8).Value = 50053.453  # H132: was 45763.523
$ws.Cells.Item(132, 9).Value = 66387  # I132: was 60588.97
$ws.Cells.Item(132, 10).Value = 4022.5454  # J132: was 3758.0833
$ws.Cells.Item(132, 11).Value = 199161  # K132: was 181766.91
$ws.Cells.Item(132, 12).Value = 12067.6362  # L132: was 11274.2499
$ws.Cells.Item(132, 13).Value = -196631  # M132: was -179236.91
$ws.Cells.Item(132, 14).Value = -17127.6362  # N132: was -16334.2499
$ws.Cells.Item(134, 8).Value = 1465.75  # H134: was 1326.5454
$ws.Cells.Item(134, 9).Value = 698.75  # I134: was 699.125
$ws.Cells.Item(134, 10).Value = 2999.75  # J134: was 2999.6667
$ws.Cells.Item(134, 11).Value = 2096.25  # K134: was 2097.375
$ws.Cells.Item(134, 12).Value = 8999.25  # L134: was 8999.000100000001
$ws.Cells.Item(134, 13).Value = 438.75  # M134: was 437.625
$ws.Cells.Item(134, 14).Value = -14069.25  # N134: was -14069.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 999.75  # H19: was 2749.75
$ws.Cells.Item(19, 10).Value = 1500  # J19: was 5000
$ws.Cells.Item(19, 12).Value = 4500  # L19: was 15000
$ws.Cells.Item(19, 14).Value = -4848  # N19: was -15348
$ws.Cells.Item(113, 8).Value = 2479.0435  # H113: was 2546.2727
$ws.Cells.Item(113, 10).Value = 3098.0625  # J113: was 3237.9333
$ws.Cells.Item(113, 12).Value = 9294.1875  # L113: was 9713.7999
$ws.Cells.Item(113, 14).Value = -13634.1875  # N113: was -14053.7999
$ws.Cells.Item(139, 8).Value = 3019.353  # H139: was 2032.2963
$ws.Cells.Item(139, 9).Value = 2866.3845  # I139: was 1866.9166
$ws.Cells.Item(139, 10).Value = 3516.5  # J139: was 3355.3333
$ws.Cells.Item(139, 11).Value = 8599.1535  # K139: was 5600.7498
$ws.Cells.Item(139, 12).Value = 10549.5  # L139: was 10065.9999
$ws.Cells.Item(139, 13).Value = -3459.1535  # M139: was -460.7497999999996
$ws.Cells.Item(139, 14).Value = -20829.5  # N139: was -20345.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 289950.28  # H70: was 254333.88
$ws.Cells.Item(70, 9).Value = 503726  # I70: was 403940.8
$ws.Cells.Item(70, 10).Value = 4916  # J70: was 4989
$ws.Cells.Item(70, 11).Value = 503726  # K70: was 403940.8
$ws.Cells.Item(70, 12).Value = 4916  # L70: was 4989
$ws.Cells.Item(70, 13).Value = -503456  # M70: was -403670.8
$ws.Cells.Item(70, 14).Value = -5456  # N70: was -5529
$ws.Cells.Item(73, 8).Value = 289950.28  # H73: was 254333.88
$ws.Cells.Item(73, 9).Value = 503726  # I73: was 403940.8
$ws.Cells.Item(73, 10).Value = 4916  # J73: was 4989
$ws.Cells.Item(73, 11).Value = 503726  # K73: was 403940.8
$ws.Cells.Item(73, 12).Value = 4916  # L73: was 4989
$ws.Cells.Item(73, 13).Value = -502790  # M73: was -403004.8
$ws.Cells.Item(73, 14).Value = -6788  # N73: was -6861
$ws.Cells.Item(97, 8).Value = 1973.3529  # H97: was 2027.1212
$ws.Cells.Item(97, 9).Value = 1471.3043  # I97: was 1529.1364
$ws.Cells.Item(97, 11).Value = 1471.3043  # K97: was 1529.1364
$ws.Cells.Item(97, 13).Value = -975.3043  # M97: was -1033.1364
$ws.Cells.Item(132, 8).Value = 2636.25  # H132: was 2450
$ws.Cells.Item(132, 9).Value = 2766.6667  # I132: was 2450
$ws.Cells.Item(132, 10).Value = 2245  # J132: was 0
$ws.Cells.Item(132, 11).Value = 8300.000100000001  # K132: was 7350
$ws.Cells.Item(132, 12).Value = 6735  # L132: was 0
$ws.Cells.Item(132, 13).Value = -5770.000100000001  # M132: was -4820
$ws.Cells.Item(132, 14).Value = -11795  # N132: new cell
$ws.Cells.Item(141, 8).Value = 119750  # H141: was 119999.5
$ws.Cells.Item(141, 10).Value = 119750  # J141: was 119999.5
$ws.Cells.Item(141, 12).Value = 119750  # L141: was 119999.5
$ws.Cells.Item(141, 14).Value = -130110  # N141: was -130359.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1215.5834  # H46: was 1128.2593
$ws.Cells.Item(46, 9).Value = 678.25  # I46: was 645.8261
$ws.Cells.Item(46, 11).Value = 678.25  # K46: was 645.8261
$ws.Cells.Item(46, 13).Value = -490.25  # M46: was -457.8261
$ws.Cells.Item(61, 8).Value = 3181.2144  # H61: was 3210.4285
$ws.Cells.Item(61, 9).Value = 2338  # I61: was 2750.4285
$ws.Cells.Item(61, 10).Value = 4699  # J61: was 3670.4285
$ws.Cells.Item(61, 11).Value = 2338  # K61: was 2750.4285
$ws.Cells.Item(61, 12).Value = 4699  # L61: was 3670.4285
$ws.Cells.Item(61, 13).Value = -2136  # M61: was -2548.4285
$ws.Cells.Item(61, 14).Value = -5103  # N61: was -4074.4285
$ws.Cells.Item(113, 8).Value = 3181.2144  # H113: was 3210.4285
$ws.Cells.Item(113, 9).Value = 2338  # I113: was 2750.4285
$ws.Cells.Item(113, 10).Value = 4699  # J113: was 3670.4285
$ws.Cells.Item(113, 11).Value = 2338  # K113: was 2750.4285
$ws.Cells.Item(113, 12).Value = 4699  # L113: was 3670.4285
$ws.Cells.Item(113, 13).Value = -168  # M113: was -580.4285
$ws.Cells.Item(113, 14).Value = -9039  # N113: was -8010.4285
$ws.Cells.Item(132, 8).Value = 12821.394  # H132: was 11959.739
$ws.Cells.Item(132, 9).Value = 11053.071  # I132: was 10022
$ws.Cells.Item(132, 11).Value = 33159.213  # K132: was 30066
$ws.Cells.Item(132, 13).Value = -30629.213  # M132: was -27536

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2538.5144  # H132: was 2328.641
$ws.Cells.Item(132, 9).Value = 3256.0833  # I132: was 2786.724
$ws.Cells.Item(132, 10).Value = 972.9091  # J132: was 1000.2
$ws.Cells.Item(132, 11).Value = 9768.249899999999  # K132: was 8360.172
$ws.Cells.Item(132, 12).Value = 2918.7273  # L132: was 3000.6
$ws.Cells.Item(132, 13).Value = -7238.249899999999  # M132: was -5830.172
$ws.Cells.Item(132, 14).Value = -7978.7273  # N132: was -8060.6
$ws.Cells.Item(136, 8).Value = 3344.3462  # H136: was 3523.625
$ws.Cells.Item(136, 9).Value = 1083.3334  # I136: was 1104.5714
$ws.Cells.Item(136, 10).Value = 6427.5454  # J136: was 6910.3
$ws.Cells.Item(136, 11).Value = 3250.0002  # K136: was 3313.7142
$ws.Cells.Item(136, 12).Value = 19282.6362  # L136: was 20730.9
$ws.Cells.Item(136, 13).Value = -700.0001999999999  # M136: was -763.7142000000003
$ws.Cells.Item(136, 14).Value = -24382.6362  # N136: was -25830.9
$ws.Cells.Item(140, 8).Value = 79569.836  # H140: was 87484.39999999999
$ws.Cells.Item(140, 10).Value = 79569.836  # J140: was 87484.39999999999
$ws.Cells.Item(140, 12).Value = 79569.836  # L140: was 87484.39999999999
$ws.Cells.Item(140, 14).Value = -89929.836  # N140: was -97844.39999999999
$ws.Cells.Item(141, 8).Value = 118191  # H141: was 109188.8
$ws.Cells.Item(141, 10).Value = 144998.33  # J141: was 129994.664
$ws.Cells.Item(141, 12).Value = 144998.33  # L141: was 129994.664
$ws.Cells.Item(141, 14).Value = -155358.33  # N141: was -140354.664
